$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header / matchup / prediction data (replaces the old 13th-Feb table
# with a refreshed set of matchups + Ballgorithm/ESPN accuracy rows).
$rows = @(
    @("NBA, Tuesday 13th Feb 2024", "Ballgorithm", "ESPN"),
    @("Atlanta Hawks (24-30) vs Charlotte Hornets (12-41)", "Atlanta Hawks (46.43%)", "Atlanta Hawks (76.0%)"),
    @("New York Knicks (33-21) vs Orlando Magic (29-25)", "New York Knicks (70.37%)", "Orlando Magic (53.6%)"),
    @("Miami Heat (29-25) vs Philadelphia 76ers (32-21)", "Philadelphia 76ers (62.96%)", "Philadelphia 76ers (68.1%)"),
    @("Brooklyn Nets (21-32) vs Boston Celtics (42-12)", "Boston Celtics (89.29%)", "Boston Celtics (83.7%)"),
    @("Chicago Bulls (26-28) vs Cleveland Cavaliers (35-17)", "Cleveland Cavaliers (66.67%)", "Cleveland Cavaliers (85.7%)"),
    @("Indiana Pacers (30-26) vs Toronto Raptors (19-35)", "Indiana Pacers (60.71%)", "Indiana Pacers (50.0%)"),
    @("Houston Rockets (24-29) vs Memphis Grizzlies (18-36)", "Houston Rockets (67.86%)", "Houston Rockets (58.0%)"),
    @("Washington Wizards (9-44) vs New Orleans Pelicans (32-22)", "New Orleans Pelicans (60.00%)", "New Orleans Pelicans (88.2%)"),
    @("San Antonio Spurs (11-43) vs Dallas Mavericks (31-23)", "Dallas Mavericks (55.17%)", "Dallas Mavericks (81.0%)"),
    @("Sacramento Kings (30-23) vs Denver Nuggets (36-18)", "Denver Nuggets (84.00%)", "Denver Nuggets (68.9%)"),
    @("Detroit Pistons (8-45) vs Phoenix Suns (32-22)", "Phoenix Suns (60.71%)", "Phoenix Suns (86.7%)"),
    @("Los Angeles Lakers (30-26) vs Utah Jazz (26-28)", "Los Angeles Lakers (68.97%)", "Utah Jazz (65.9%)"),
    @("Los Angeles Clippers (35-17) vs Golden State Warriors (26-25)", "Los Angeles Clippers (76.92%)", "Golden State Warriors (50.7%)")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $rowVals = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
}

$ws.Range("C15").Select()
